$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1578.5454
$ws.Range("I4").Value = 1552.7142
$ws.Range("K4").Value = 1552.7142
$ws.Range("M4").Value = -1438.7142
$ws.Range("H92").Value = 1203.5555
$ws.Range("I92").Value = 1499
$ws.Range("K92").Value = 1499
$ws.Range("M92").Value = -251
$ws.Range("H107").Value = 345.83334
$ws.Range("I107").Value = 345.83334
$ws.Range("K107").Value = 345.83334
$ws.Range("M107").Value = 1574.16666
$ws.Range("H132").Value = 1867.7693
$ws.Range("I132").Value = 1867.7693
$ws.Range("K132").Value = 5603.3079
$ws.Range("M132").Value = -3073.3079

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1123
$ws.Range("H61").Value = 1204.125
$ws.Range("I61").Value = 1204.125
$ws.Range("K61").Value = 1204.125
$ws.Range("M61").Value = -992.125
$ws.Range("H74").Value = 3669.7
$ws.Range("I74").Value = 3077.4443
$ws.Range("K74").Value = 3077.4443
$ws.Range("M74").Value = -2203.4443
$ws.Range("H77").Value = 3669.7
$ws.Range("I77").Value = 3077.4443
$ws.Range("K77").Value = 15387.2215
$ws.Range("M77").Value = -11019.2215
$ws.Range("H80").Value = 133596.6
$ws.Range("I80").Value = 68998
$ws.Range("K80").Value = 68998
$ws.Range("M80").Value = -68000
$ws.Range("H83").Value = 133596.6
$ws.Range("I83").Value = 68998
$ws.Range("K83").Value = 206994
$ws.Range("M83").Value = -202002
$ws.Range("H136").Value = 1204.125
$ws.Range("I136").Value = 1204.125
$ws.Range("K136").Value = 3612.375
$ws.Range("M136").Value = -1062.375

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 800.25
$ws.Range("I22").Value = 174.5
$ws.Range("K22").Value = 174.5
$ws.Range("M22").Value = -1.5
$ws.Range("H75").Value = 33960
$ws.Range("I75").Value = 33960
$ws.Range("K75").Value = 33960
$ws.Range("M75").Value = -33024
$ws.Range("H76").Value = 92985.875
$ws.Range("J76").Value = 92985.875
$ws.Range("L76").Value = 92985.875
$ws.Range("N76").Value = -93615.875
$ws.Range("H78").Value = 33960
$ws.Range("I78").Value = 33960
$ws.Range("K78").Value = 101880
$ws.Range("M78").Value = -97200
$ws.Range("H79").Value = 92985.875
$ws.Range("J79").Value = 92985.875
$ws.Range("L79").Value = 92985.875
$ws.Range("N79").Value = -95169.875
$ws.Range("H86").Value = 1733
$ws.Range("I86").Value = 1755.3334
$ws.Range("J86").Value = 1666
$ws.Range("K86").Value = 1755.3334
$ws.Range("L86").Value = 1666
$ws.Range("M86").Value = -632.3334
$ws.Range("N86").Value = -3912
$ws.Range("H88").Value = 21383.857
$ws.Range("I88").Value = 8498
$ws.Range("J88").Value = 23531.5
$ws.Range("K88").Value = 8498
$ws.Range("L88").Value = 23531.5
$ws.Range("M88").Value = -8092
$ws.Range("N88").Value = -24343.5
$ws.Range("H89").Value = 1733
$ws.Range("I89").Value = 1755.3334
$ws.Range("J89").Value = 1666
$ws.Range("K89").Value = 8776.666999999999
$ws.Range("L89").Value = 8330
$ws.Range("M89").Value = -3160.666999999999
$ws.Range("N89").Value = -19562
$ws.Range("H91").Value = 21383.857
$ws.Range("I91").Value = 8498
$ws.Range("J91").Value = 23531.5
$ws.Range("K91").Value = 8498
$ws.Range("L91").Value = 23531.5
$ws.Range("M91").Value = -7094
$ws.Range("N91").Value = -26339.5
$ws.Range("H107").Value = 471
$ws.Range("I107").Value = 446.75
$ws.Range("J107").Value = 519.5
$ws.Range("K107").Value = 446.75
$ws.Range("L107").Value = 519.5
$ws.Range("M107").Value = 1473.25
$ws.Range("N107").Value = -4359.5
$ws.Range("H134").Value = 2070.7368
$ws.Range("I134").Value = 2002.75
$ws.Range("K134").Value = 6008.25
$ws.Range("M134").Value = -3473.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 947.5
$ws.Range("I15").Value = 947.5
$ws.Range("K15").Value = 947.5
$ws.Range("M15").Value = -777.5
$ws.Range("H81").Value = 78999.664
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 78999.664
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 78999.664
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -80995.664
$ws.Range("H84").Value = 78999.664
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 78999.664
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 236998.992
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -246982.992
$ws.Range("H99").Value = 2966.4
$ws.Range("I99").Value = 2966.4
$ws.Range("K99").Value = 2966.4
$ws.Range("M99").Value = -1468.4
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H122").Value = 1939.2106
$ws.Range("I122").Value = 1765.3125
$ws.Range("J122").Value = 2866.6667
$ws.Range("K122").Value = 5295.9375
$ws.Range("L122").Value = 8600.000100000001
$ws.Range("M122").Value = -2845.9375
$ws.Range("N122").Value = -13500.0001
$ws.Range("H126").Value = 2966.4
$ws.Range("I126").Value = 2966.4
$ws.Range("K126").Value = 8899.200000000001
$ws.Range("M126").Value = -6429.200000000001
$ws.Range("H132").Value = 3302.8823
$ws.Range("I132").Value = 3353.7144
$ws.Range("K132").Value = 10061.1432
$ws.Range("M132").Value = -7531.143199999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13750895
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 300
$ws.Range("N4").Value = -524
$ws.Range("H33").Value = 444
$ws.Range("I33").Value = 40.333332
$ws.Range("K33").Value = 241.999992
$ws.Range("M33").Value = 41.00000800000001
$ws.Range("H122").Value = 442.57144
$ws.Range("I122").Value = 319.8
$ws.Range("K122").Value = 2878.2
$ws.Range("M122").Value = -428.2000000000003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2353.6
$ws.Range("I102").Value = 2439.25
$ws.Range("J102").Value = 2011
$ws.Range("K102").Value = 2439.25
$ws.Range("L102").Value = 2011
$ws.Range("M102").Value = -817.25
$ws.Range("N102").Value = -5255
$ws.Range("H132").Value = 2249.375
$ws.Range("I132").Value = 2249.375
$ws.Range("K132").Value = 6748.125
$ws.Range("M132").Value = -4218.125

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 20751
$ws.Range("J103").Value = 20751
$ws.Range("L103").Value = 20751
$ws.Range("N103").Value = -23095
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
$ws.Range("H132").Value = 7480.0835
$ws.Range("I132").Value = 8740.223
$ws.Range("K132").Value = 26220.669
$ws.Range("M132").Value = -23690.669

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1977
$ws.Range("I122").Value = 1520
$ws.Range("J122").Value = 2548.25
$ws.Range("K122").Value = 4560
$ws.Range("L122").Value = 7644.75
$ws.Range("M122").Value = -2110
$ws.Range("N122").Value = -12544.75
$ws.Range("H132").Value = 1884.7222
$ws.Range("J132").Value = 1429.6666
$ws.Range("L132").Value = 4288.9998
$ws.Range("N132").Value = -9348.9998

